# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price cells that look like plain decimals (e.g. "65.79") are prefixed with a
# leading apostrophe so Excel stores them as text instead of re-parsing them as
# numbers (which would strip significant trailing zeros / introduce float
# rounding noise). Values that already aren't valid numbers (e.g. "26.881.19",
# thousand-grouped prices) and the percentage/name/link cells are assigned as
# plain strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.881.19'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '1.669.25'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '''215.66'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('E6').Value = '  +4.41%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = '''20.27'
$ws.Range('E10').Value = '  +3.81%  '
$ws.Range('E11').Value = '  +3.84%  '
$ws.Range('D12').Value = '1.904.55'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '1.697.37'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '''65.79'
$ws.Range('E16').Value = '  +1.25%  '
$ws.Range('D17').Value = '26.901.90'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '''231.96'
$ws.Range('E18').Value = '  -4.01%  '
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('D22').Value = '''4.47'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = '''2.22'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '''9.20'
$ws.Range('E24').Value = '  -0.65%  '
$ws.Range('D25').Value = '''145.58'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = '''0.116'
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '''7.13'
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('D28').Value = '''15.91'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = '''1.18'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').Value = '''3.34'
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('D33').Value = '1.465.00'
$ws.Range('E33').Value = '  -3.80%  '
$ws.Range('D34').Value = '''3.17'
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('E35').Value = '  +4.12%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').Value = '''0.899'
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').Value = '''5.81'
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '''0.973'
$ws.Range('E43').Value = '  +6.43%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''65.74'
$ws.Range('E44').Value = '  +1.30%  '
$ws.Range('D45').Value = '1.813.90'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('D47').Value = '''90.40'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('E49').Value = '  +2.15%  '
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').Value = '''7.59'
$ws.Range('E51').Value = '  +0.89%  '
